$wb = $excel.ActiveWorkbook

# --- ALC row block (hunk idx 0) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 61.636364
$ws.Range("I8").Value = 69.333336
$ws.Range("J8").Value = 27
$ws.Range("K8").Value = 208.000008
$ws.Range("L8").Value = 81
$ws.Range("M8").Value = -69.00000800000001
$ws.Range("N8").Value = -359

# --- ALC row block (hunk idx 1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1209.7273
$ws.Range("I28").Value = 862.125
$ws.Range("K28").Value = 862.125
$ws.Range("M28").Value = -377.125

# --- ALC row block (hunk idx 2) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 353.72726
$ws.Range("J55").Value = 383.66666
$ws.Range("L55").Value = 383.66666
$ws.Range("N55").Value = -811.66666

# --- ALC row block (hunk idx 3) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 7878.1763
$ws.Range("I64").Value = 4468
$ws.Range("K64").Value = 4468
$ws.Range("M64").Value = -4220

# --- ALC row block (hunk idx 4) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 7878.1763
$ws.Range("I67").Value = 4468
$ws.Range("K67").Value = 4468
$ws.Range("M67").Value = -3610

# --- ALC row block (hunk idx 5) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1688.7273
$ws.Range("I80").Value = 1002.7778
$ws.Range("J80").Value = 2163.6155
$ws.Range("K80").Value = 3008.3334
$ws.Range("L80").Value = 6490.8465
$ws.Range("M80").Value = -2010.3334
$ws.Range("N80").Value = -8486.8465

# --- ALC row block (hunk idx 6) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1688.7273
$ws.Range("I83").Value = 1002.7778
$ws.Range("J83").Value = 2163.6155
$ws.Range("K83").Value = 9025.0002
$ws.Range("L83").Value = 19472.5395
$ws.Range("M83").Value = -4033.0002
$ws.Range("N83").Value = -29456.5395

# --- ALC row block (hunk idx 7) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 6053.4546
$ws.Range("J86").Value = 7248.8335
$ws.Range("L86").Value = 7248.8335
$ws.Range("N86").Value = -9494.833500000001

# --- ALC row block (hunk idx 8) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 6053.4546
$ws.Range("J89").Value = 7248.8335
$ws.Range("L89").Value = 36244.1675
$ws.Range("N89").Value = -47476.1675

# --- ALC row block (hunk idx 9) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1838.3334
$ws.Range("J103").Value = 2128.889
$ws.Range("L103").Value = 6386.667
$ws.Range("N103").Value = -7558.667

# --- ALC row block (hunk idx 10) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3674.6
$ws.Range("I106").Value = 3674.6
$ws.Range("K106").Value = 3674.6
$ws.Range("M106").Value = -3043.6

# --- ALC row block (hunk idx 11) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1606.52
$ws.Range("I132").Value = 1338.15
$ws.Range("K132").Value = 4014.45
$ws.Range("M132").Value = -1484.45

# --- ALC row block (hunk idx 12) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1667.05
$ws.Range("I137").Value = 1657.8462
$ws.Range("J137").Value = 1684.1428
$ws.Range("K137").Value = 4973.5386
$ws.Range("L137").Value = 5052.428400000001
$ws.Range("M137").Value = -2423.5386
$ws.Range("N137").Value = -10152.4284

# --- ALC row block (hunk idx 13) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 786.1429000000001
$ws.Range("I141").Value = 775.8889
$ws.Range("K141").Value = 2327.6667
$ws.Range("M141").Value = 2852.3333

# --- ARM row block (hunk idx 14) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1297.125
$ws.Range("I74").Value = 1297.125
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1297.125
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -423.125
$ws.Range("N74").Value = $null

# --- ARM row block (hunk idx 15) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1297.125
$ws.Range("I77").Value = 1297.125
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 6485.625
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -2117.625
$ws.Range("N77").Value = $null

# --- ARM row block (hunk idx 16) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2794
$ws.Range("I102").Value = 2794
$ws.Range("K102").Value = 2794
$ws.Range("M102").Value = -1172

# --- ARM row block (hunk idx 17) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1149.1666
$ws.Range("I132").Value = 299.66666
$ws.Range("K132").Value = 898.9999799999999
$ws.Range("M132").Value = 1631.00002

# --- BSM row block (hunk idx 18) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 32918
$ws.Range("I26").Value = 32918
$ws.Range("K26").Value = 32918
$ws.Range("M26").Value = -32626

# --- BSM row block (hunk idx 19) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11247.5
$ws.Range("I86").Value = 9495
$ws.Range("K86").Value = 9495
$ws.Range("M86").Value = -8372

# --- BSM row block (hunk idx 20) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 11247.5
$ws.Range("I89").Value = 9495
$ws.Range("K89").Value = 47475
$ws.Range("M89").Value = -41859

# --- BSM row block (hunk idx 21) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3666.3333
$ws.Range("I107").Value = 3666.3333
$ws.Range("K107").Value = 3666.3333
$ws.Range("M107").Value = -1746.3333

# --- BSM row block (hunk idx 22) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1996.4117
$ws.Range("I134").Value = 2039.9375
$ws.Range("K134").Value = 6119.8125
$ws.Range("M134").Value = -3584.8125

# --- CRP row block (hunk idx 23) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 546.35297
$ws.Range("I22").Value = 555.1111
$ws.Range("J22").Value = 536.5
$ws.Range("K22").Value = 555.1111
$ws.Range("L22").Value = 536.5
$ws.Range("M22").Value = -205.1111
$ws.Range("N22").Value = -1236.5

# --- CRP row block (hunk idx 24) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2399
$ws.Range("J62").Value = 2399
$ws.Range("L62").Value = 2399
$ws.Range("N62").Value = -3647

# --- CRP row block (hunk idx 25) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2399
$ws.Range("J65").Value = 2399
$ws.Range("L65").Value = 11995
$ws.Range("N65").Value = -18235

# --- CRP row block (hunk idx 26) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1672.6086
$ws.Range("I134").Value = 1242.8823
$ws.Range("K134").Value = 3728.6469
$ws.Range("M134").Value = -1193.6469

# --- CUL row block (hunk idx 27) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 70.25
$ws.Range("I13").Value = 62
$ws.Range("J13").Value = 95
$ws.Range("K13").Value = 186
$ws.Range("L13").Value = 285
$ws.Range("M13").Value = -18
$ws.Range("N13").Value = -621

# --- CUL row block (hunk idx 28) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 840.7143
$ws.Range("I113").Value = 770.3333
$ws.Range("K113").Value = 2310.9999
$ws.Range("M113").Value = -140.9998999999998

# --- CUL row block (hunk idx 29) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 3045.3333
$ws.Range("I121").Value = 3696.3333
$ws.Range("K121").Value = 11088.9999
$ws.Range("M121").Value = -9778.999899999999

# --- CUL row block (hunk idx 30) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 277659.5
$ws.Range("I128").Value = 277659.5
$ws.Range("K128").Value = 832978.5
$ws.Range("M128").Value = -827998.5

# --- CUL row block (hunk idx 31) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1713.3334
$ws.Range("J129").Value = 3196.6667
$ws.Range("L129").Value = 9590.000100000001
$ws.Range("N129").Value = -19590.0001

# --- CUL row block (hunk idx 32) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 14750
$ws.Range("J130").Value = 15000
$ws.Range("L130").Value = 45000
$ws.Range("N130").Value = -55040

# --- CUL row block (hunk idx 33) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 953.2105
$ws.Range("I131").Value = 873.5
$ws.Range("J131").Value = 990
$ws.Range("K131").Value = 2620.5
$ws.Range("L131").Value = 2970
$ws.Range("M131").Value = 2419.5
$ws.Range("N131").Value = -13050

# --- GSM row block (hunk idx 34) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3877.45
$ws.Range("I122").Value = 3871.5557
$ws.Range("K122").Value = 11614.6671
$ws.Range("M122").Value = -9164.667099999999

# --- LTW row block (hunk idx 35) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1300
$ws.Range("I7").Value = 1100
$ws.Range("K7").Value = 1100
$ws.Range("M7").Value = -988

# --- LTW row block (hunk idx 36) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 401.66666
$ws.Range("I55").Value = 82
$ws.Range("J55").Value = 2000
$ws.Range("K55").Value = 82
$ws.Range("L55").Value = 2000
$ws.Range("M55").Value = 91
$ws.Range("N55").Value = -2346

# --- LTW row block (hunk idx 37) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2130.25
$ws.Range("I122").Value = 2166.4285
$ws.Range("K122").Value = 6499.2855
$ws.Range("M122").Value = -4049.2855

# --- LTW row block (hunk idx 38) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1300
$ws.Range("I126").Value = 1100
$ws.Range("K126").Value = 3300
$ws.Range("M126").Value = -830

# --- WVR row block (hunk idx 39) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3665
$ws.Range("I122").Value = 3665
$ws.Range("K122").Value = 10995
$ws.Range("M122").Value = -8545

# --- WVR row block (hunk idx 40) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1885.125
$ws.Range("I136").Value = 1992.8462
$ws.Range("J136").Value = 1418.3334
$ws.Range("K136").Value = 5978.5386
$ws.Range("L136").Value = 4255.0002
$ws.Range("M136").Value = -3428.5386
$ws.Range("N136").Value = -9355.0002
